# Edit script generated from OOXML diff for the FFXIV Leve profits workbook.
# For each affected sheet/row, updates the currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H,I,J,K,L,M,N) to the new values from the commit. Cells that the diff removes
# entirely (no replacement <v>) are cleared instead of set to 0/blank text.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H3").Value = 45000
$ws.Range("J3").Value = 45000
$ws.Range("L3").Value = 45000
$ws.Range("N3").Value = -45228
$ws.Range("H28").Value = 1104.2
$ws.Range("I28").Value = 1054.5
$ws.Range("J28").Value = 1800
$ws.Range("K28").Value = 1054.5
$ws.Range("L28").Value = 1800
$ws.Range("M28").Value = -569.5
$ws.Range("N28").Value = -2770
$ws.Range("H32").Value = 467.8
$ws.Range("J32").Value = 384.75
$ws.Range("L32").Value = 384.75
$ws.Range("N32").Value = -1036.75
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("N95").Value = 0
$ws.Range("H102").Value = 45000
$ws.Range("J102").Value = 45000
$ws.Range("L102").Value = 45000
$ws.Range("N102").Value = -51490
$ws.Range("H113").Value = 13341.529
$ws.Range("J113").Value = 17083.834
$ws.Range("L113").Value = 17083.834
$ws.Range("N113").Value = -23591.834
$ws.Range("H137").Value = 2182.532
$ws.Range("I137").Value = 1492.8684
$ws.Range("J137").Value = 5094.4443
$ws.Range("K137").Value = 4478.6052
$ws.Range("L137").Value = 15283.3329
$ws.Range("M137").Value = -1928.6052
$ws.Range("N137").Value = -20383.3329
$ws.Range("H138").Value = 2680.7913
$ws.Range("I138").Value = 1749
$ws.Range("J138").Value = 2910.5479
$ws.Range("K138").Value = 5247
$ws.Range("L138").Value = 8731.643700000001
$ws.Range("M138").Value = -107
$ws.Range("N138").Value = -19011.6437

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 608.7143
$ws.Range("I2").Value = 541.3333
$ws.Range("K2").Value = 541.3333
$ws.Range("M2").Value = -428.3333
$ws.Range("H116").Value = 608.7143
$ws.Range("I116").Value = 541.3333
$ws.Range("K116").Value = 541.3333
$ws.Range("M116").Value = 1752.6667
$ws.Range("H122").Value = 3052.8
$ws.Range("I122").Value = 1828.5714
$ws.Range("J122").Value = 4124
$ws.Range("K122").Value = 5485.7142
$ws.Range("L122").Value = 12372
$ws.Range("M122").Value = -3035.7142
$ws.Range("N122").Value = -17272
$ws.Range("H133").Value = 43880.25
$ws.Range("J133").Value = 43880.25
$ws.Range("L133").Value = 43880.25
$ws.Range("N133").Value = -48940.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H3").Value = 608.7143
$ws.Range("I3").Value = 541.3333
$ws.Range("K3").Value = 541.3333
$ws.Range("M3").Value = -427.3333
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492
$ws.Range("H99").Value = 2725.6333
$ws.Range("I99").Value = 1049.3125
$ws.Range("J99").Value = 4641.4287
$ws.Range("K99").Value = 1049.3125
$ws.Range("L99").Value = 4641.4287
$ws.Range("M99").Value = 448.6875
$ws.Range("N99").Value = -7637.4287
$ws.Range("H134").Value = 2771.3584
$ws.Range("I134").Value = 1573.5238
$ws.Range("J134").Value = 7344.909
$ws.Range("K134").Value = 4720.5714
$ws.Range("L134").Value = 22034.727
$ws.Range("M134").Value = -2185.5714
$ws.Range("N134").Value = -27104.727

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H9").Value = 32420
$ws.Range("J9").Value = 32420
$ws.Range("L9").Value = 32420
$ws.Range("N9").Value = -32756
$ws.Range("H16").Value = 6946762.5
$ws.Range("I16").Value = 22223240
$ws.Range("J16").Value = 2909.0908
$ws.Range("K16").Value = 22223240
$ws.Range("L16").Value = 2909.0908
$ws.Range("M16").Value = -22222953
$ws.Range("N16").Value = -3483.0908
$ws.Range("H58").Value = 1939.638
$ws.Range("I58").Value = 1710.5088
$ws.Range("J58").Value = 15000
$ws.Range("K58").Value = 1710.5088
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -1507.5088
$ws.Range("N58").Value = -15406
$ws.Range("H74").Value = 34858.23
$ws.Range("J74").Value = 34858.23
$ws.Range("L74").Value = 34858.23
$ws.Range("N74").Value = -36606.23
$ws.Range("H77").Value = 34858.23
$ws.Range("J77").Value = 34858.23
$ws.Range("L77").Value = 104574.69
$ws.Range("N77").Value = -113310.69
$ws.Range("H96").Value = 15000
$ws.Range("J96").Value = 15000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -20492
$ws.Range("H113").Value = 6946762.5
$ws.Range("I113").Value = 22223240
$ws.Range("J113").Value = 2909.0908
$ws.Range("K113").Value = 22223240
$ws.Range("L113").Value = 2909.0908
$ws.Range("M113").Value = -22221070
$ws.Range("N113").Value = -7249.0908
$ws.Range("H122").Value = 1957.0344
$ws.Range("I122").Value = 1427.7084
$ws.Range("J122").Value = 4497.8
$ws.Range("K122").Value = 4283.1252
$ws.Range("L122").Value = 13493.4
$ws.Range("M122").Value = -1833.1252
$ws.Range("N122").Value = -18393.4
$ws.Range("H132").Value = 2713.92
$ws.Range("I132").Value = 1645.3125
$ws.Range("J132").Value = 4613.6665
$ws.Range("K132").Value = 4935.9375
$ws.Range("L132").Value = 13840.9995
$ws.Range("M132").Value = -2405.9375
$ws.Range("N132").Value = -18900.9995
$ws.Range("H136").Value = 1939.638
$ws.Range("I136").Value = 1710.5088
$ws.Range("J136").Value = 15000
$ws.Range("K136").Value = 5131.526400000001
$ws.Range("L136").Value = 45000
$ws.Range("M136").Value = -2581.526400000001
$ws.Range("N136").Value = -50100

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H80").Value = 5322.5
$ws.Range("I80").Value = 1980
$ws.Range("J80").Value = 5800
$ws.Range("K80").Value = 5940
$ws.Range("L80").Value = 17400
$ws.Range("M80").Value = -5004
$ws.Range("N80").Value = -19272
$ws.Range("H83").Value = 5322.5
$ws.Range("I83").Value = 1980
$ws.Range("J83").Value = 5800
$ws.Range("K83").Value = 17820
$ws.Range("L83").Value = 52200
$ws.Range("M83").Value = -13140
$ws.Range("N83").Value = -61560

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H123").Value = 15262.75
$ws.Range("J123").Value = 15262.75
$ws.Range("L123").Value = 15262.75
$ws.Range("N123").Value = -20162.75
$ws.Range("H126").Value = 3993.0989
$ws.Range("I126").Value = 2864.04
$ws.Range("J126").Value = 5370
$ws.Range("K126").Value = 8592.119999999999
$ws.Range("L126").Value = 16110
$ws.Range("M126").Value = -6122.119999999999
$ws.Range("N126").Value = -21050
$ws.Range("H132").Value = 2917.6667
$ws.Range("I132").Value = 1660
$ws.Range("J132").Value = 3057.4075
$ws.Range("K132").Value = 4980
$ws.Range("L132").Value = 9172.2225
$ws.Range("M132").Value = -2450
$ws.Range("N132").Value = -14232.2225

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 6498.75
$ws.Range("I7").Value = 4700
$ws.Range("J7").Value = 7398.125
$ws.Range("K7").Value = 4700
$ws.Range("L7").Value = 7398.125
$ws.Range("M7").Value = -4588
$ws.Range("N7").Value = -7622.125
$ws.Range("H16").Value = 875
$ws.Range("I16").Value = 875
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 875
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -705
$ws.Range("H40").Value = 5797.8184
$ws.Range("I40").Value = 5316
$ws.Range("K40").Value = 5316
$ws.Range("M40").Value = -5180
$ws.Range("H61").Value = 2438.0715
$ws.Range("I61").Value = 2344.4167
$ws.Range("K61").Value = 2344.4167
$ws.Range("M61").Value = -2142.4167
$ws.Range("H113").Value = 2438.0715
$ws.Range("I113").Value = 2344.4167
$ws.Range("K113").Value = 2344.4167
$ws.Range("M113").Value = -174.4167000000002
$ws.Range("H122").Value = 3612.8044
$ws.Range("I122").Value = 2948.3513
$ws.Range("J122").Value = 6344.4443
$ws.Range("K122").Value = 8845.053899999999
$ws.Range("L122").Value = 19033.3329
$ws.Range("M122").Value = -6395.053899999999
$ws.Range("N122").Value = -23933.3329
$ws.Range("H126").Value = 6498.75
$ws.Range("I126").Value = 4700
$ws.Range("J126").Value = 7398.125
$ws.Range("K126").Value = 14100
$ws.Range("L126").Value = 22194.375
$ws.Range("M126").Value = -11630
$ws.Range("N126").Value = -27134.375
$ws.Range("H132").Value = 5514.6387
$ws.Range("I132").Value = 1791.2667
$ws.Range("J132").Value = 8174.1904
$ws.Range("K132").Value = 5373.800099999999
$ws.Range("L132").Value = 24522.5712
$ws.Range("M132").Value = -2843.800099999999
$ws.Range("N132").Value = -29582.5712
$ws.Range("H133").Value = 49837.727
$ws.Range("J133").Value = 49837.727
$ws.Range("L133").Value = 49837.727
$ws.Range("N133").Value = -54897.727
$ws.Range("H139").Value = 43206.555
$ws.Range("J139").Value = 43607.375
$ws.Range("L139").Value = 43607.375
$ws.Range("N139").Value = -53887.375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").Value = 0
$ws.Range("H82").Value = 45650
$ws.Range("J82").Value = 45650
$ws.Range("L82").Value = 45650
$ws.Range("N82").Value = -46416
$ws.Range("H85").Value = 45650
$ws.Range("J85").Value = 45650
$ws.Range("L85").Value = 45650
$ws.Range("N85").Value = -48302
$ws.Range("H126").Value = 636159.8
$ws.Range("I126").Value = 3985.25
$ws.Range("K126").Value = 11955.75
$ws.Range("M126").Value = -9485.75
$ws.Range("H132").Value = 8134461.5
$ws.Range("I132").Value = 5322.773
$ws.Range("J132").Value = 17547148
$ws.Range("K132").Value = 15968.319
$ws.Range("L132").Value = 52641444
$ws.Range("M132").Value = -13438.319
$ws.Range("N132").Value = -52646504
$ws.Range("H136").Value = 4377.577
$ws.Range("I136").Value = 1370.75
$ws.Range("J136").Value = 6954.857
$ws.Range("K136").Value = 4112.25
$ws.Range("L136").Value = 20864.571
$ws.Range("M136").Value = -1562.25
$ws.Range("N136").Value = -25964.571
$ws.Range("H139").Value = 40489.656
$ws.Range("J139").Value = 40872.223
$ws.Range("L139").Value = 40872.223
$ws.Range("N139").Value = -51152.223

